$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 353, pushing the existing rows 353..376 down to 355..378.
$ws.Rows.Item(353).Insert()
$ws.Rows.Item(353).Insert()

# --- New row 353 ---
$ws.Cells.Item(353, 1).Value = 10
$ws.Cells.Item(353, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(353, 3).Value = "La Araucanía"
$ws.Cells.Item(353, 4).Value = 44746
$ws.Cells.Item(353, 5).Value = 9
$ws.Cells.Item(353, 6).Value = "Fruta"
$ws.Cells.Item(353, 7).Value = 100108
$ws.Cells.Item(353, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(353, 9).Value = 100108002
$ws.Cells.Item(353, 10).Value = "Mango"
$ws.Cells.Item(353, 11).Value = "Sin especificar"
$ws.Cells.Item(353, 12).Value = "Primera"
$ws.Cells.Item(353, 13).Value = 1500
$ws.Cells.Item(353, 14).Value = 10000
$ws.Cells.Item(353, 15).Value = 10000
$ws.Cells.Item(353, 16).Value = 10000
$ws.Cells.Item(353, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(353, 18).Value = "Brasil"
$ws.Cells.Item(353, 19).Value = 2500
$ws.Cells.Item(353, 20).Value = 4

# --- New row 354 ---
$ws.Cells.Item(354, 1).Value = 10
$ws.Cells.Item(354, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(354, 3).Value = "La Araucanía"
$ws.Cells.Item(354, 4).Value = 44746
$ws.Cells.Item(354, 5).Value = 9
$ws.Cells.Item(354, 6).Value = "Fruta"
$ws.Cells.Item(354, 7).Value = 100108
$ws.Cells.Item(354, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(354, 9).Value = 100108002
$ws.Cells.Item(354, 10).Value = "Mango"
$ws.Cells.Item(354, 11).Value = "Sin especificar"
$ws.Cells.Item(354, 12).Value = "Segunda"
$ws.Cells.Item(354, 13).Value = 500
$ws.Cells.Item(354, 14).Value = 8000
$ws.Cells.Item(354, 15).Value = 8000
$ws.Cells.Item(354, 16).Value = 8000
$ws.Cells.Item(354, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(354, 18).Value = "Brasil"
$ws.Cells.Item(354, 19).Value = 2000
$ws.Cells.Item(354, 20).Value = 4
